$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as plain TEXT (matching the workbook's
# existing convention of storing every data cell as an inline/shared string,
# even when the text looks numeric). Temporarily marking the cell as Text
# before the write stops Excel's auto-number-detection; ClearFormats()
# afterwards drops the temporary Text number-format again so no stray style
# is left attached to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Insert a new column at G ("P_B"), pushing the existing G (EPS(TTM)) and H
# (Div. Yield) columns one to the right, to H and I respectively. Insert()
# shifts cell content/styles automatically, including the bold header style.
$ws.Columns("G:G").Insert()

# New header cell for the inserted column.
$ws.Range("G1").Value = "P_B"

# Updated data row - ticker swapped from Coca-Cola to Philip Morris, plus
# refreshed figures and the new P/B value in column G.
Set-TextValue $ws.Range("A2") "Philip Morris International Inc"
Set-TextValue $ws.Range("B2") "`$155.16"
Set-TextValue $ws.Range("C2") "+0.15(0.10%) 1D"
Set-TextValue $ws.Range("D2") "`$241.25B"
Set-TextValue $ws.Range("F2") "34.2"
Set-TextValue $ws.Range("G2") "-20.5"
Set-TextValue $ws.Range("H2") "1.32"
Set-TextValue $ws.Range("I2") "3.48%"
